$d = $word.ActiveDocument

$newSentence = " One often has difficulty with this at the beginning, so we provide a number of examples and also give you some exercises to do yourself."

# Avoid inserting twice if the script is ever re-applied to an already-edited document.
$already = $d.Content.Find.Execute("One often has difficulty with this at the beginning", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $already) {

    # Locate the target sentence ending in "execution of an algorithm." so the
    # insertion point is resolved dynamically rather than via a hard-coded offset.
    $find = $d.Content.Find
    $found = $find.Execute("execution of an algorithm.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Could not find target sentence 'execution of an algorithm.'"
    }

    # Collapse the found range to its end point (right after the final period,
    # i.e. right before the _GoBack bookmark that follows it) and insert the
    # new sentence there as a brand-new run.
    $insertPoint = $find.Parent.Duplicate
    $insertPoint.Collapse(0)
    $insertPoint.InsertAfter($newSentence)

    # The new text is visually identical in formatting to the preceding run
    # (Times New Roman, 10pt), so the engine would normally merge it into
    # that same run. Explicitly (re)stamping the font forces it to remain a
    # distinct <w:r> with its own <w:rPr>, matching the target edit exactly.
    $insertPoint.Font.Name = "Times New Roman"
    $insertPoint.Font.NameBi = "Times New Roman"
}
